# "Chapters 8 to 10"
#
# - Remove the run of slides covering the "operational definitions /
#   experiment 1 / experiment 2" detour (old slides 10-16), leaving the
#   "For Monday Jan 22" wrap-up slide as the new slide 10.
# - Retitle the opening slide for the new class date.
# - Tweak two lines of text on the (now-last) wrap-up slide for the new
#   exam date / wording.

$p = $ppt.ActivePresentation

# --- Drop the 7 slides between "Within-participant design concerns /
#     Questions" (slide 9) and the "For Monday Jan 22" wrap-up slide.
#     They keep sliding into position 10 as each one is removed.
for ($i = 0; $i -lt 7; $i++) {
    $p.Slides.Item(10).Delete()
}

# --- Slide 1 title: "205 Oct 7, Class 8" -> "205 Jan 19, Class 8"
$slide1 = $p.Slides.Item(1)
$titleRange = $slide1.Shapes.Item(1).TextFrame.TextRange
$titleRange.Runs(1).Text = "205 Jan 19, Class 8"

# --- Wrap-up slide (now slide 10) body text tweaks.
$lastSlide = $p.Slides.Item(10)
$bodyRange = $lastSlide.Shapes.Item(2).TextFrame.TextRange

$examPara = $bodyRange.Paragraphs(7)
$examPara.Runs(1).Text = "Wednesday Jan 24: Exam 1"

$studyAidsPara = $bodyRange.Paragraphs(8)
$studyAidsPara.Runs(1).Text = "Prior exams are posted to Canvas as study aids"
